$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (was Paris / ECCMID 2020 / Apr-20) to the new
# Montpellier / EEID 2020 / Jun-20 conference entry.
$ws.Range("A6").Value = "Montpellier"
$ws.Range("B6").Value = "EEID 2020"
$ws.Range("C6").Value = "Jun-20"
$ws.Range("D6").Value = "Will attend"
$ws.Range("E6").Value = 43.6
$ws.Range("F6").Value = 3.8833

# Update the active selection to match the saved workbook state.
$ws.Range("C5").Select()
